$d = $word.ActiveDocument

# The document contains five "<id>...</id>" tags, each currently split
# across three runs: "<id>" (Courier New / 7f6000 / sz18), the bare id
# value (black, default font), and "</id>" (Courier New / 7f6000 / sz18).
# Collapse each trio into a single run "<id>VALUE</id>" that keeps the
# formatting of the opening "<id>" run.

$searchStart = 0
for ($i = 1; $i -le 5; $i++) {

    # Locate the opening "<id>" run starting after the previous match.
    $openRange = $d.Range($searchStart, $d.Content.End)
    $openFound = $openRange.Find.Execute("<id>", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

    # Locate the matching closing "</id>" run.
    $closeRange = $d.Range($openRange.End, $d.Content.End)
    $closeFound = $closeRange.Find.Execute("</id>", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

    $openStart = $openRange.Start
    $openEnd   = $openRange.End
    $fullEnd   = $closeRange.End

    # Full text spanning "<id>" through "</id>" (all 3 runs).
    $fullText = $d.Range($openStart, $fullEnd).Text

    # Re-typing the opening run's text with the full merged string makes
    # it adopt the opening run's own formatting (Courier New/7f6000/sz18).
    $openRange.Text = $fullText
    $newOpenEnd = $openStart + $fullText.Length

    # Remove the now-duplicated leftover text of the old middle/closing runs.
    $leftoverLength = $fullEnd - $openEnd
    $leftover = $d.Range($newOpenEnd, $newOpenEnd + $leftoverLength)
    $leftover.Delete()

    $searchStart = $newOpenEnd
}
